$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.Range("C3:D5")
$rng.VerticalAlignment = -4108
